$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Common Word"
$ws.Range("B1").Value = "Total Frequency"
$ws.Range("C1").Value = "Websites"

# Data rows
$ws.Range("A2").Value = "A3"
$ws.Range("B2").Value = 71
$ws.Range("C2").Value = "http://www.yyy.at/ (54), https://www.xxx.at/ (6), https://www.zzz.at/ (11)"

$ws.Range("A3").Value = "B3"
$ws.Range("B3").Value = 68
$ws.Range("C3").Value = "http://www.yyy.at/ (52), https://www.xxx.at/ (6), https://www.zzz.at/ (10)"

$ws.Range("A4").Value = "C2"
$ws.Range("B4").Value = 55
$ws.Range("C4").Value = "http://www.yyy.at/ (50), https://www.xxx.at/ (5)"

$ws.Range("A5").Value = "D2"
$ws.Range("B5").Value = 13
$ws.Range("C5").Value = "https://www.xxx.at/ (5), https://www.zzz.at/ (8)"

# Build the header style on a single cell first (bold font, thin box border,
# centered horizontally, top-aligned vertically), then fan it out to the
# rest of the header row by copying formats only - this reuses one style
# record instead of minting a new one per cell.
$h = $ws.Range("A1")
$h.Font.Bold = $true
$h.HorizontalAlignment = -4108
$h.VerticalAlignment = -4160
$h.Borders.LineStyle = 1

$h.Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
